$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the full previous range (A1:B60) before rewriting the shrunk table (A1:B49)
$ws.Range("A1:B60").ClearContents()

$data = @(
    @('Cluster Name', 'Activecases'),
    @('3398 BlueCross Elly Kay Mordialloc', 29),
    @('3601 Baptcare Westhaven community', 13),
    @('3647 Aurrum Aged Care Reservoir', 11),
    @('3653 Fronditha Thalpori St Albans Aged Care', 20),
    @('3975 Aurrum Aged Care Brunswick West', 13),
    @('4257 BlueCross The Gables Camberwell', 16),
    @('4295 Hope Aged Care Sunshine West', 14),
    @('4314 Estia Health Ardeer', 17),
    @('44095 Myrniong Primary School Myrniong', 12),
    @('44404 Castlemaine North Primary SchoolCastlemaine', 65),
    @('44507 Inverloch Primary School Inverloch', 10),
    @('44593 Torquay P-6 College Torquay', 23),
    @('44622 Grey Street Primary School Traralgon', 13),
    @('44642 Irymple South Primary School Irymple South', 14),
    @('4479 Whittlesea Lodge Whittlesea', 16),
    @('44855 Beverley Hills Primary School Doncaster East', 11),
    @('44893 Greenhills Primary School Greensborough', 13),
    @('44960 Thomastown West Primary School', 14),
    @('45013 Gladstone Views Primary School', 10),
    @('45168 Ranfurly Primary School', 37),
    @('45305 Lockington Consolidated SchoolLockington', 43),
    @('4574 Village Glen Aged Care Residences Mornington', 10),
    @('45757 Saint Joseph''s Primary School Warragul', 12),
    @('45958 Ave Maria College Aberfeldie Workplace', 11),
    @('50412 Geelong Grammar School Bostock House Campus Newtown', 11),
    @('52390 Our Lady of the Way Catholic PrimarySchool Wallan', 49),
    @('52694 Pakenham Primary School Pakenham', 20),
    @('52777 Mirripoa Primary School Mount Duneed School Camp', 26),
    @('Cardinia Waters Retirement Village Pakenham', 16),
    @('Confirmed Omicron Sircuit Bar Fitzroy', 20),
    @('Confirmed Omicron Variant The Peel HotelCollingwood', 18),
    @('Flight QF10 LHR to MEL via DRW 13 Dec', 13),
    @('Goodstart Early Learning Preston', 10),
    @('Green Gables Lodge Warburton', 27),
    @('Greendale Hotel Greendale', 23),
    @('JBS Australia Brooklyn', 35),
    @('Jackson''s Hotel Age of Love Event Toorak', 12),
    @('PGL Camp Rumbug Foster North', 15),
    @('Social Gathering 11 Dec Windsor', 23),
    @('Social Gathering 11 December Fitzroy', 10),
    @('St Brigid''s Parish Primary School Mordialloc', 14),
    @('St Clare''s Primary School Officer', 10),
    @('St Vincents Hospital Melbourne Emergency Department Fitzroy', 11),
    @('StarTrack- Melbourne Tullamarine', 20),
    @('Sunny Ridge Strawberry Farm Main Ridge', 20),
    @('The George Lounge St Kilda', 12),
    @('Thomastown West Primary School Camp Doxa''s Malmsbury', 12),
    @('Warbuton Lodge Warbuton', 15)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
